$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap columns B and C: header text and observed data values
$ws.Range("B1").Value = "Conc  [mg/l]"
$ws.Range("C1").Value = "Error  [mg/l]"

$ws.Range("B2").Value = 0.80699998140335083
$ws.Range("C2").Value = 0.18999999761581421

$ws.Range("B3").Value = 0.97999995946884155
$ws.Range("C3").Value = 0.35799998044967651

$ws.Range("B4").Value = 1.3300000429153442
$ws.Range("C4").Value = 0.38999998569488525

$ws.Range("B5").Value = 1.4199999570846558
$ws.Range("C5").Value = 0.37800002098083496
